$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 46075 -> 46076) for every
# data row (rows 2 through 426). Increment each cell by one day.
for ($r = 2; $r -le 426; $r++) {
    $ws.Cells.Item($r, 3).Value = 46076
}
